# "Agrego TS, BO y FF a secuencia"
# Adds two new rows (9 and 10) to the sequence table describing the new
# Training-Strategy / BO (Bayesian Optimization) steps.
#
# NOTE on write order: the values below are written in the same order the
# original author's Excel session created them (C9, B9, B10, C10, D10, F9,
# E10, F10) so the shared-string table ends up populated in that exact
# sequence, matching the target workbook byte-for-byte in that respect.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 - "Training strategy" step
$ws.Range("C9").Value = "Todos los campos del paso anterior"
$ws.Range("B9").Value = "Training strategy. Separar datasets: bo, train y test"
$ws.Range("D9").Value = 931
$ws.Range("E9").Value = "exp/FE9252/dataset.csv.gz"

# Row 10 - "BO" step (started before F9/E10/F10 below, same as original author)
$ws.Range("B10").Value = "BO"
$ws.Range("C10").Value = "Dataset del paso anterior"
$ws.Range("D10").Value = "941_HT"

# F9 holds a multi-line destination list (wrapped cell, taller row)
$ws.Range("F9").Value = "exp/TS9310/dataset_future.csv.gz #futuro`nexp/TS9310/dataset_train_final.csv.gz #para entrenar modelo final`nexp/TS9310/dataset_training.csv.gz #para hacer bo"

$ws.Range("E10").Value = "exp/TS9310/dataset_training.csv.gz"
$ws.Range("F10").Value = "exp/HT9410/dataset_training.csv.gz"

# Formatting: F9 wraps text and row 9 grows to fit the 3 lines
$ws.Range("F9").WrapText = $true
$ws.Rows.Item(9).RowHeight = 46.5

# Column F is widened to a fixed 60 characters (no longer auto-fit).
# ColumnWidth goes through a character<->stored-width conversion that adds
# 5/6 to whatever is assigned, so we back that offset out here to land on
# exactly 60 in the saved file.
$ws.Columns.Item(6).ColumnWidth = 60 - (5/6)

# Selection ends on C9, matching where the author's cursor landed
$ws.Range("C9").Select() | Out-Null
